# Update cryptocurrency Price (column D) and Volume(1h) (column E) values.
# A leading apostrophe forces Excel to treat the value as text (matching the
# original inline-string cells), and ClearFormats() removes the quote-prefix
# style marker that Excel would otherwise attach to the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.285.79"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'  +2.55%  "
$ws.Range("E2").ClearFormats()
$ws.Range("D3").Value = "'1.870.53"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'  +1.22%  "
$ws.Range("E3").ClearFormats()
$ws.Range("E4").Value = "'  -0.17%  "
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = "'339.91"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'  +2.12%  "
$ws.Range("E5").ClearFormats()
$ws.Range("E6").Value = "'  -0.22%  "
$ws.Range("E6").ClearFormats()
$ws.Range("D7").Value = "'0.4711"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "'  +1.51%  "
$ws.Range("E7").ClearFormats()
$ws.Range("D8").Value = "'0.3937"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "'  +2.04%  "
$ws.Range("E8").ClearFormats()
$ws.Range("D9").Value = "'47.28"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "'  +2.74%  "
$ws.Range("E9").ClearFormats()
$ws.Range("D10").Value = "'0.07997"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "'  +0.97%  "
$ws.Range("E10").ClearFormats()
$ws.Range("D11").Value = "'1.009"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "'  +1.44%  "
$ws.Range("E11").ClearFormats()
$ws.Range("E12").Value = "'  +1.92%  "
$ws.Range("E12").ClearFormats()
$ws.Range("D13").Value = "'1.881.85"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "'  +1.04%  "
$ws.Range("E13").ClearFormats()
$ws.Range("D14").Value = "'6.009"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'  +1.57%  "
$ws.Range("E14").ClearFormats()
$ws.Range("D15").Value = "'7.287"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "'  +2.60%  "
$ws.Range("E15").ClearFormats()
$ws.Range("D16").Value = "'91.24"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'  +2.74%  "
$ws.Range("E16").ClearFormats()
$ws.Range("D17").Value = "'1.002"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "'  -0.21%  "
$ws.Range("E17").ClearFormats()
$ws.Range("D18").Value = "'0.00001040"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "'  +0.48%  "
$ws.Range("E18").ClearFormats()
$ws.Range("D19").Value = "'0.06596"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "'  -0.91%  "
$ws.Range("E19").ClearFormats()
$ws.Range("D20").Value = "'17.68"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "'  +3.59%  "
$ws.Range("E20").ClearFormats()
$ws.Range("D21").Value = "'1.002"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "'  -0.14%  "
$ws.Range("E21").ClearFormats()
$ws.Range("D22").Value = "'28.297.18"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "'  +2.59%  "
$ws.Range("E22").ClearFormats()
$ws.Range("D23").Value = "'5.460"
$ws.Range("D23").ClearFormats()
$ws.Range("E24").Value = "'  +1.30%  "
$ws.Range("E24").ClearFormats()
$ws.Range("D25").Value = "'2.288"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "'  -0.52%  "
$ws.Range("E25").ClearFormats()
$ws.Range("D26").Value = "'2.098.24"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "'  +1.24%  "
$ws.Range("E26").ClearFormats()
$ws.Range("E27").Value = "'  +1.16%  "
$ws.Range("E27").ClearFormats()
$ws.Range("E28").Value = "'  +1.53%  "
$ws.Range("E28").ClearFormats()
$ws.Range("D29").Value = "'2.145"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "'  +2.24%  "
$ws.Range("E29").ClearFormats()
$ws.Range("D30").Value = "'5.504"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "'  +1.92%  "
$ws.Range("E30").ClearFormats()
$ws.Range("D31").Value = "'120.47"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "'  +0.57%  "
$ws.Range("E31").ClearFormats()
$ws.Range("D32").Value = "'0.9774"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "'  +0.34%  "
$ws.Range("E32").ClearFormats()
$ws.Range("D33").Value = "'0.09505"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "'  +1.17%  "
$ws.Range("E33").ClearFormats()
$ws.Range("D34").Value = "'3.595"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "'  +0.63%  "
$ws.Range("E34").ClearFormats()
$ws.Range("D35").Value = "'1.380"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "'  +2.56%  "
$ws.Range("E35").ClearFormats()
$ws.Range("D36").Value = "'5.358"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "'  +1.38%  "
$ws.Range("E36").ClearFormats()
$ws.Range("D37").Value = "'0.02278"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "'  +2.51%  "
$ws.Range("E37").ClearFormats()
$ws.Range("D38").Value = "'0.06097"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "'  +1.33%  "
$ws.Range("E38").ClearFormats()
$ws.Range("D39").Value = "'8.456"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "'  +1.83%  "
$ws.Range("E39").ClearFormats()
$ws.Range("E40").Value = "'  -0.07%  "
$ws.Range("E40").ClearFormats()
$ws.Range("E41").Value = "'  +1.36%  "
$ws.Range("E41").ClearFormats()
$ws.Range("D42").Value = "'1.000"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "'  -0.04%  "
$ws.Range("E42").ClearFormats()
$ws.Range("E43").Value = "'  +1.07%  "
$ws.Range("E43").ClearFormats()
$ws.Range("D44").Value = "'10.40"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "'  +1.01%  "
$ws.Range("E44").ClearFormats()
$ws.Range("D45").Value = "'1.308"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "'  +4.82%  "
$ws.Range("E45").ClearFormats()
$ws.Range("D46").Value = "'0.5611"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "'  +0.43%  "
$ws.Range("E46").ClearFormats()
$ws.Range("D47").Value = "'12.14"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "'  -0.45%  "
$ws.Range("E47").ClearFormats()
$ws.Range("D48").Value = "'1.968"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "'  +3.69%  "
$ws.Range("E48").ClearFormats()
$ws.Range("E49").Value = "'  +3.35%  "
$ws.Range("E49").ClearFormats()
$ws.Range("E50").Value = "'  +0.51%  "
$ws.Range("E50").ClearFormats()
$ws.Range("D51").Value = "'2.025"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "'  +13.54%  "
$ws.Range("E51").ClearFormats()
